# "ultima versión de datos"
# Fill in the two newly-reported "indicador 6 / Ranking" rows (Yucatán and
# Chiapas, year 2022) that were still missing their value in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
[void]$ws.Activate()

$ws.Range("E318").Value = 2
$ws.Range("E319").Value = 1

# Restore the default zoom (the author's last save had it reset to 100%)
# and leave the selection where the author left it before saving.
$win = $excel.ActiveWindow
$win.Zoom = 100
[void]$ws.Range("A846").Select()
[void]$ws.Range("C326").Select()
